$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
foreach ($p in @(@("H13",2582.8), @("I13",1816.3334), @("J13",2911.2856), @("K13",1816.3334), @("L13",2911.2856), @("M13",-1647.3334), @("N13",-3249.2856))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H32",50000), @("I32",0), @("J32",50000), @("K32",0), @("L32",50000), @("N32",-50652))) { $ws.Range($p[0]).Value = $p[1] }
$ws.Range("M32").ClearContents()

foreach ($p in @(@("H41",982.625), @("I41",640.53845), @("K41",640.53845), @("M41",-200.53845))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H70",2978.6667), @("I70",1594), @("J70",4363.3335), @("K70",4782), @("L70",13090.0005), @("M70",-4512), @("N70",-13630.0005))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H73",2978.6667), @("I73",1594), @("J73",4363.3335), @("K73",4782), @("L73",13090.0005), @("M73",-3846), @("N73",-14962.0005))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H81",81194.5), @("J81",81194.5), @("L81",81194.5), @("N81",-83190.5))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H84",81194.5), @("J84",81194.5), @("L84",243583.5), @("N84",-253567.5))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H107",3291.3635), @("I107",3075.25), @("J107",3414.8572), @("K107",3075.25), @("L107",3414.8572), @("M107",-1155.25), @("N107",-7254.8572))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H138",1084794.2), @("I138",588.9167), @("J138",1518476.4), @("K138",1766.7501), @("L138",4555429.199999999), @("M138",3373.2499), @("N138",-4565709.199999999))) { $ws.Range($p[0]).Value = $p[1] }

$ws = $wb.Worksheets.Item("ARM")
foreach ($p in @(@("H45",1684.65), @("I45",911.44446), @("J45",2317.2727), @("K45",911.44446), @("L45",2317.2727), @("M45",-534.44446), @("N45",-3071.2727))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H61",45553588), @("I61",125001530), @("J61",154763.42), @("K61",125001530), @("L61",154763.42), @("M61",-125001318), @("N61",-155187.42))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H74",10007245), @("I74",14706715), @("J74",20869.25), @("K74",14706715), @("L74",20869.25), @("M74",-14705841), @("N74",-22617.25))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H77",10007245), @("I77",14706715), @("J77",20869.25), @("K77",73533575), @("L77",104346.25), @("M77",-73529207), @("N77",-113082.25))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H132",5842.1787), @("I132",3128.5), @("J132",12626.375), @("K132",9385.5), @("L132",37879.125), @("M132",-6855.5), @("N132",-42939.125))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H136",45553588), @("I136",125001530), @("J136",154763.42), @("K136",375004590), @("L136",464290.26), @("M136",-375002040), @("N136",-469390.26))) { $ws.Range($p[0]).Value = $p[1] }

$ws = $wb.Worksheets.Item("BSM")
foreach ($p in @(@("H99",2072), @("I99",1720.9), @("J99",2949.75), @("K99",1720.9), @("L99",2949.75), @("M99",-222.9000000000001), @("N99",-5945.75))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H134",40056.31), @("I134",1504.3182), @("K134",4512.9546), @("M134",-1977.9546))) { $ws.Range($p[0]).Value = $p[1] }

$ws = $wb.Worksheets.Item("CRP")
foreach ($p in @(@("H16",1356), @("I16",990), @("K16",990), @("M16",-703))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H28",112999.5), @("J28",112999.5), @("L28",112999.5), @("N28",-113489.5))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H31",758687.4399999999), @("I31",11364.267), @("J31",1459303), @("K31",11364.267), @("L31",1459303), @("M31",-11069.267), @("N31",-1459893))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H34",758687.4399999999), @("I34",11364.267), @("J34",1459303), @("K34",11364.267), @("L34",1459303), @("M34",-11162.267), @("N34",-1459707))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H58",2486.1177), @("I58",2441.6667), @("J58",2592.8), @("K58",2441.6667), @("L58",2592.8), @("M58",-2238.6667), @("N58",-2998.8))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H96",1479.6), @("J96",1349.75), @("L96",1349.75), @("N96",-6841.75))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H113",1356), @("I113",990), @("K113",990), @("M113",1180))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H132",3976.3125), @("I132",4044.7144), @("J132",3497.5), @("K132",12134.1432), @("L132",10492.5), @("M132",-9604.143199999999), @("N132",-15552.5))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H134",403273.1), @("I134",527780.5600000001), @("J134",8999.333000000001), @("K134",1583341.68), @("L134",26997.999), @("M134",-1580806.68), @("N134",-32067.999))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H136",2486.1177), @("I136",2441.6667), @("J136",2592.8), @("K136",7325.000100000001), @("L136",7778.400000000001), @("M136",-4775.000100000001), @("N136",-12878.4))) { $ws.Range($p[0]).Value = $p[1] }

$ws = $wb.Worksheets.Item("CUL")
foreach ($p in @(@("H4",13379091), @("I4",9376288), @("K4",28128864), @("M4",-28128752))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H56",10757.462), @("I56",10757.462), @("K56",10757.462), @("M56",-10227.462))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H112",4577.222), @("I112",2456.8572), @("J112",11998.5), @("K112",7370.571599999999), @("L112",35995.5), @("M112",-6262.571599999999), @("N112",-38211.5))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H137",3131.7058), @("I137",2621.1667), @("J137",4357), @("K137",7863.500100000001), @("L137",13071), @("M137",-2763.500100000001), @("N137",-23271))) { $ws.Range($p[0]).Value = $p[1] }

$ws = $wb.Worksheets.Item("GSM")
foreach ($p in @(@("H70",4800.8), @("J70",9997), @("L70",9997), @("N70",-10537))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H73",4800.8), @("J73",9997), @("L73",9997), @("N73",-11869))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H92",21371.166), @("J92",21371.166), @("L92",21371.166), @("N92",-25115.166))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H102",2293.4167), @("I102",2306.4546), @("J102",2150), @("K102",2306.4546), @("L102",2150), @("M102",-684.4546), @("N102",-5394))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H126",2313.4443), @("I126",2260.5715), @("K126",6781.7145), @("M126",-4311.7145))) { $ws.Range($p[0]).Value = $p[1] }

$ws = $wb.Worksheets.Item("LTW")
foreach ($p in @(@("H7",53962.75), @("I7",3112.818), @("K7",3112.818), @("M7",-3000.818))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H45",4820.5), @("I45",4820.5), @("K45",4820.5), @("M45",-4413.5))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H46",3203.0303), @("I46",2600.4), @("K46",2600.4), @("M46",-2412.4))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H97",99996.664), @("J97",99996.664), @("L97",99996.664), @("N97",-101978.664))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H106",44500), @("I106",59000), @("J106",39666.668), @("K106",59000), @("L106",39666.668), @("M106",-57738), @("N106",-42190.668))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H126",53962.75), @("I126",3112.818), @("K126",9338.454000000002), @("M126",-6868.454000000002))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H132",277385.06), @("I132",280191.62), @("J132",252126), @("K132",840574.86), @("L132",756378), @("M132",-838044.86), @("N132",-761438))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H136",80227.92999999999), @("I136",47256.773), @("J136",225301), @("K136",141770.319), @("L136",675903), @("M136",-139220.319), @("N136",-681003))) { $ws.Range($p[0]).Value = $p[1] }

$ws = $wb.Worksheets.Item("WVR")
foreach ($p in @(@("H92",62200), @("J92",62200), @("L92",62200), @("N92",-67192))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H97",0), @("J97",0), @("L97",0))) { $ws.Range($p[0]).Value = $p[1] }
$ws.Range("N97").ClearContents()

foreach ($p in @(@("H107",31251070), @("I107",45455796), @("J107",677.8), @("K107",136367388), @("L107",2033.4), @("M107",-136365468), @("N107",-5873.4))) { $ws.Range($p[0]).Value = $p[1] }

foreach ($p in @(@("H136",1808.7693), @("I136",1608.762), @("J136",2648.8), @("K136",4826.286), @("L136",7946.400000000001), @("M136",-2276.286), @("N136",-13046.4))) { $ws.Range($p[0]).Value = $p[1] }
